$wb = $excel.ActiveWorkbook

# Insert a new worksheet "Arcs2" right after "Arcs" (it becomes the active sheet,
# matching the tabSelected flag moving off of "Arcs" and onto the new sheet).
$arcsSheet = $wb.Worksheets.Item("Arcs")
$arcs2 = $wb.Worksheets.Add($null, $arcsSheet)
$arcs2.Name = "Arcs2"

# Header row
$arcs2.Range("A1").Value = "Arc"
$arcs2.Range("B1").Value = "From"
$arcs2.Range("C1").Value = "To"
$arcs2.Range("D1").Value = "Cost"
$arcs2.Range("E1").Value = "Capacity"

$data = @(
    @(0,  "ANC", "CVG", 10, 5),
    @(1,  "ANC", "LAX", 7,  5),
    @(2,  "ANC", "ICN", 15, 2),
    @(3,  "ANC", "PVG", 13, 3),
    @(4,  "ANC", "HKG", 18, 4),
    @(5,  "LAX", "ICN", 12, 8),
    @(6,  "CVG", "HKG", 20, 4),
    @(7,  "CVG", "LEJ", 17, 2),
    @(8,  "CVG", "EMA", 15, 2),
    @(9,  "JFK", "LEJ", 19, 8),
    @(10, "JFK", "EMA", 17, 6),
    @(11, "EMA", "LEJ", 3,  10),
    @(12, "EMA", "BAH", 10, 10),
    @(13, "LEJ", "BRU", 4,  8),
    @(14, "LEJ", "LOS", 9,  12),
    @(15, "LEJ", "BAH", 10, 5),
    @(16, "LEJ", "DEL", 13, 10),
    @(17, "LEJ", "ICN", 22, 5),
    @(18, "LEJ", "HKG", 25, 3),
    @(19, "LEJ", "DXB", 8,  5),
    @(20, "BRU", "BAH", 12, 5),
    @(21, "BAH", "PVG", 11, 14),
    @(22, "BAH", "HKG", 13, 4),
    @(23, "BAH", "BKK", 12, 8),
    @(24, "BAH", "SIN", 9,  10),
    @(25, "DXB", "HKG", 10, 5),
    @(26, "DEL", "SIN", 6,  8),
    @(27, "BKK", "SIN", 3,  7),
    @(28, "ICN", "HKG", 5,  10),
    @(29, "ICN", "PVG", 4,  5)
)

$row = 2
foreach ($r in $data) {
    $arcs2.Cells.Item($row, 1).Value = $r[0]
    $arcs2.Cells.Item($row, 2).Value = $r[1]
    $arcs2.Cells.Item($row, 3).Value = $r[2]
    $arcs2.Cells.Item($row, 4).Value = $r[3]
    $arcs2.Cells.Item($row, 5).Value = $r[4]
    $row++
}

# The whole used range is center-aligned in the original sheet style.
$arcs2.Range("A1:E31").HorizontalAlignment = -4108

# Row 14 (arc 12, EMA-BAH) has its From/To cells rendered in explicit black font.
$arcs2.Range("B14:C14").Font.Color = 0

# Match the view state captured in the saved file.
$arcs2.Range("C31").Select()
$av = $excel.ActiveWindow
$av.ScrollRow = 14
